$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C12").Value = "BABY NADIF NASAL SPRAY 50 ML"
$ws.Range("H12").Value = "0:0"
$ws.Range("L12").Value = "1"
$ws.Range("N12").Value = "218.00"
$ws.Range("P12").Value = "218.0000"
$ws.Range("Q12").Value = "1:0"
$ws.Range("C13").Value = "CETAL 250MG/5ML 60ML SUSP"
$ws.Range("H13").Value = "13:0"
$ws.Range("L13").Value = "1"
$ws.Range("N13").Value = "31.00"
$ws.Range("P13").Value = "31.0000"
$ws.Range("Q13").Value = "1:0"
$ws.Range("C14").Value = "DECANCIT S.R 20 F.C.TAB"
$ws.Range("H14").Value = "0:1"
$ws.Range("L14").Value = "1"
$ws.Range("N14").Value = "42.00"
$ws.Range("P14").Value = "21.0000"
$ws.Range("Q14").Value = "0:1"
$ws.Range("C15").Value = "DEPOVIT B12-1000MCG/ML 5 I.M. AMP"
$ws.Range("H15").Value = "0:4"
$ws.Range("L15").Value = "1"
$ws.Range("N15").Value = "85.00"
$ws.Range("P15").Value = "17.0000"
$ws.Range("Q15").Value = "0:1"
$ws.Range("C16").Value = "DEXAMETHASONE-AMRIYA 8MG/2ML 3 AMP."
$ws.Range("H16").Value = "0:0"
$ws.Range("L16").Value = "1"
$ws.Range("N16").Value = "36.00"
$ws.Range("P16").Value = "23.7600"
$ws.Range("Q16").Value = "0:2"
$ws.Range("C17").Value = "DICLAC 75 ID 30 TAB"
$ws.Range("H17").Value = "1:1"
$ws.Range("L17").Value = "0"
$ws.Range("N17").Value = "135.00"
$ws.Range("P17").Value = "44.5500"
$ws.Range("Q17").Value = "0:1"
$ws.Range("C18").Value = "DIMRA 20 F.C.TAB."
$ws.Range("H18").Value = "0:1"
$ws.Range("L18").Value = "1"
$ws.Range("N18").Value = "70.00"
$ws.Range("P18").Value = "35.0000"
$ws.Range("Q18").Value = "0:1"
$ws.Range("C19").Value = "ERASTAPEX PLUS 20MG/12.5MG 30 TAB"
$ws.Range("H19").Value = "0:0"
$ws.Range("L19").Value = "1"
$ws.Range("N19").Value = "78.00"
$ws.Range("P19").Value = "25.7400"
$ws.Range("Q19").Value = "0:1"
$ws.Range("C20").Value = "EZAMOL-C 20 TAB."
$ws.Range("H20").Value = "1:0"
$ws.Range("L20").Value = "1"
$ws.Range("N20").Value = "24.00"
$ws.Range("P20").Value = "12.0000"
$ws.Range("Q20").Value = "0:1"
$ws.Range("C21").Value = "FAROVIGA 100MG 12 F.C.TAB."
$ws.Range("H21").Value = "2:1"
$ws.Range("L21").Value = "1"
$ws.Range("N21").Value = "108.00"
$ws.Range("P21").Value = "35.6400"
$ws.Range("Q21").Value = "0:4"
$ws.Range("A49").Value = "Thursday, 7 August, 2025 2:12 PM"
$ws.Range("P48").Value = 2346.8800000000001
